# Generate Report for handback
# Updates the zh-cn and de-de status sheets to reflect that the handback
# for "00be777f-..." and "2a4ef3df-..." has been completed:
#   - Status changes from "Ready for handoff" to "Handed back: in sync with en-US"
#   - Latest Target File / Latest Handback File columns (E/F) get populated
#     with links to the handoff markdown / xlf files (mirroring columns A/C)
#   - Latest Handback DateTime (G) is stamped with the handback time

$wb = $excel.ActiveWorkbook

function Update-StatusSheet($SheetName, $Row2TargetUrl, $Row2HandbackUrl, $Row2HandbackFileName, $Row3TargetUrl, $Row3HandbackUrl, $Row3HandbackFileName, $Row2DateTime, $Row3DateTime) {

    $ws = $wb.Worksheets.Item($SheetName)

    # --- Row 2 (00be777f-2241-48cc-bfc5-feff68d1127e.md) ---
    $ws.Range("B2").Value = "Handed back: in sync with en-US"

    $ws.Range("E2").Value = "00be777f-2241-48cc-bfc5-feff68d1127e.md"
    $ws.Hyperlinks.Add($ws.Range("E2"), $Row2TargetUrl, "", "", "00be777f-2241-48cc-bfc5-feff68d1127e.md") | Out-Null

    $ws.Range("F2").Value = $Row2HandbackFileName
    $ws.Hyperlinks.Add($ws.Range("F2"), $Row2HandbackUrl, "", "", $Row2HandbackFileName) | Out-Null

    $ws.Range("G2").Value = $Row2DateTime

    # --- Row 3 (2a4ef3df-2a97-4c83-90a1-9918c8ba7ad3.md) ---
    $ws.Range("B3").Value = "Handed back: in sync with en-US"

    $ws.Range("E3").Value = "2a4ef3df-2a97-4c83-90a1-9918c8ba7ad3.md"
    $ws.Hyperlinks.Add($ws.Range("E3"), $Row3TargetUrl, "", "", "2a4ef3df-2a97-4c83-90a1-9918c8ba7ad3.md") | Out-Null

    $ws.Range("F3").Value = $Row3HandbackFileName
    $ws.Hyperlinks.Add($ws.Range("F3"), $Row3HandbackUrl, "", "", $Row3HandbackFileName) | Out-Null

    $ws.Range("G3").Value = $Row3DateTime
}

# zh-cn sheet
Update-StatusSheet `
    "zh-cn" `
    "https://github.com/OpenLocalizationTest/oltest/blob/2fb0d14f85b8e0f10dd3cfd8ce9244336e3c88cb/e2e/00be777f-2241-48cc-bfc5-feff68d1127e.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/38fcde09ec66015f95defd1565e77a888d94dd6c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/00be777f-2241-48cc-bfc5-feff68d1127e.3eff22196a8fd0708f82cdd78211739061baba16.zh-cn.xlf" `
    "00be777f-2241-48cc-bfc5-feff68d1127e.3eff22196a8fd0708f82cdd78211739061baba16.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTest/oltest/blob/2fb0d14f85b8e0f10dd3cfd8ce9244336e3c88cb/e2e/2a4ef3df-2a97-4c83-90a1-9918c8ba7ad3.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/38fcde09ec66015f95defd1565e77a888d94dd6c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/2a4ef3df-2a97-4c83-90a1-9918c8ba7ad3.d437bb8ab6e675cda2bab9132481d5635525fc4b.zh-cn.xlf" `
    "2a4ef3df-2a97-4c83-90a1-9918c8ba7ad3.d437bb8ab6e675cda2bab9132481d5635525fc4b.zh-cn.xlf" `
    "2016-01-26 06:28:32" `
    "2016-01-26 06:28:32"

# de-de sheet
Update-StatusSheet `
    "de-de" `
    "https://github.com/OpenLocalizationTest/oltest/blob/2fb0d14f85b8e0f10dd3cfd8ce9244336e3c88cb/e2e/00be777f-2241-48cc-bfc5-feff68d1127e.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/659df0b75abbe1b1ac596b7b8341ad98c4e8b5fe/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/00be777f-2241-48cc-bfc5-feff68d1127e.3eff22196a8fd0708f82cdd78211739061baba16.de-de.xlf" `
    "00be777f-2241-48cc-bfc5-feff68d1127e.3eff22196a8fd0708f82cdd78211739061baba16.de-de.xlf" `
    "https://github.com/OpenLocalizationTest/oltest/blob/2fb0d14f85b8e0f10dd3cfd8ce9244336e3c88cb/e2e/2a4ef3df-2a97-4c83-90a1-9918c8ba7ad3.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/659df0b75abbe1b1ac596b7b8341ad98c4e8b5fe/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/2a4ef3df-2a97-4c83-90a1-9918c8ba7ad3.d437bb8ab6e675cda2bab9132481d5635525fc4b.de-de.xlf" `
    "2a4ef3df-2a97-4c83-90a1-9918c8ba7ad3.d437bb8ab6e675cda2bab9132481d5635525fc4b.de-de.xlf" `
    "2016-01-26 06:28:54" `
    "2016-01-26 06:28:54"

$wb.Save()
